$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as literal text even when it looks like a number
# (e.g. "580.40"), preserving the original (default) cell style afterwards.
function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

# Apply the refreshed crypto price / 1h-volume data, including the
# Bittensor <-> dogwifhat row re-ranking (rows 43-44).
$ws.Range("D2").Value = "67.062.41"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.120.97"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue "D5" "580.40"
$ws.Range("E5").Value = "  -0.10%  "
Set-TextValue "D6" "173.65"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.31%  "
Set-TextValue "D9" "6.42"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("E11").Value = "  -0.76%  "
Set-TextValue "D12" "0.0000249"
$ws.Range("E12").Value = "  -0.24%  "
Set-TextValue "D13" "37.21"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "3.637.63"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "67.041.33"
$ws.Range("E16").Value = "  +0.21%  "
Set-TextValue "D17" "7.13"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "3.120.73"
$ws.Range("E18").Value = "  +0.36%  "
Set-TextValue "D19" "16.41"
$ws.Range("E19").Value = "  +2.09%  "
Set-TextValue "D20" "491.29"
$ws.Range("E20").Value = "  +1.79%  "
Set-TextValue "D21" "7.96"
$ws.Range("E21").Value = "  +5.81%  "
Set-TextValue "D22" "0.707"
$ws.Range("E22").Value = "  -1.08%  "
Set-TextValue "D23" "84.12"
$ws.Range("E23").Value = "  +0.10%  "
Set-TextValue "D24" "13.23"
$ws.Range("E24").Value = "  +0.84%  "
Set-TextValue "D25" "2.29"
$ws.Range("E25").Value = "  -3.52%  "
Set-TextValue "D26" "10.37"
$ws.Range("E26").Value = "  +3.09%  "
$ws.Range("E27").Value = "  +0.01%  "
Set-TextValue "D28" "7.91"
$ws.Range("E28").Value = "  -0.86%  "
Set-TextValue "D29" "2.36"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("E30").Value = "  -0.49%  "
Set-TextValue "D31" "28.60"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").Value = "0.0₃0947"
$ws.Range("E33").Value = "  -6.11%  "
$ws.Range("E34").Value = "  -0.12%  "
Set-TextValue "D35" "5.88"
$ws.Range("E35").Value = "  -0.51%  "
Set-TextValue "D36" "0.976"
$ws.Range("E36").Value = "  -1.86%  "
Set-TextValue "D37" "47.46"
$ws.Range("E37").Value = "  -1.49%  "
Set-TextValue "D38" "2.06"
$ws.Range("E38").Value = "  -3.07%  "
Set-TextValue "D39" "0.311"
$ws.Range("E39").Value = "  -2.29%  "
Set-TextValue "D40" "0.124"
$ws.Range("E40").Value = "  +1.39%  "
Set-TextValue "D41" "8.53"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").Value = "2.822.61"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "2.61"
$ws.Range("E43").Value = "  -7.27%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D44" "383.40"
$ws.Range("E44").Value = "  -0.15%  "
Set-TextValue "D45" "0.0353"
$ws.Range("E45").Value = "  -2.42%  "
Set-TextValue "D46" "135.50"
$ws.Range("E46").Value = "  +0.48%  "
Set-TextValue "D48" "25.01"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("E50").Value = "  -0.80%  "
Set-TextValue "D51" "6.75"
